# Model_Iterations.xlsx — add a 6th model iteration row documenting the
# standard-scaler / lasso / ridge experiment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append the new iteration as row 7 -------------------------------------
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 0.59389000000000003
$ws.Cells.Item(7, 3).Value = "Joel Rosario, Eric Cancel, Manuel Franco, Andre Shivnarie Worrie, Reylu Gutierres"
$ws.Cells.Item(7, 4).Value = "track_id, race_number, race_date, jockey, program_number"
$ws.Cells.Item(7, 5).Value = "race_date, latitude, longitude, trakus_index, program_number, distance_id, run_up_distance, purse, post_time, odds"
$ws.Cells.Item(7, 6).Value = "Standaridizes data using standard scalar. Tests lasso and ridge regressions (multivariate linear regression still the best in terms of predicitve power). "

# Match the wrapped-text formatting used by the other data rows (columns C:F)
$ws.Range("C7:F7").WrapText = $true

# Row grew tall enough (5 wrapped lines of notes) to need an explicit height
$ws.Rows.Item(7).RowHeight = 85

# --- Update the recorded selection/view to the new active cell -------------
$ws.Range("G6").Select() | Out-Null
